$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Range("A5:K5")
$row.NumberFormat = "@"

$ws.Range("A5").Value = "Melbourne"
$ws.Range("B5").Value = "Shanghai"
$ws.Range("C5").Value = "40REHC"
$ws.Range("D5").Value = "800"
$ws.Range("E5").Value = "400"
$ws.Range("F5").Value = "500"
$ws.Range("G5").Value = "400"
$ws.Range("H5").Value = "50"
$ws.Range("I5").Value = "60"
$ws.Range("J5").Value = "Collect"
$ws.Range("K5").Value = "14 Days"

$row.Style = "Normal"
